$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for plain (unstyled) data cells in column D, used to
# strip the "quote prefix" style Excel applies when a text value is
# forced via a leading apostrophe (keeps cell formatting identical to
# the original inline-string cells).
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "'73.090.90"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").Value = "'4.002.56"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'592.30"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +9.55%  "

$ws.Range("D6").Value = "'160.78"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +7.48%  "

$ws.Range("D7").Value = "'0.688"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.753"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +1.88%  "

$ws.Range("D10").Value = "'0.169"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +2.18%  "

$ws.Range("D11").Value = "'54.21"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -3.65%  "

$ws.Range("D12").Value = "'0.0000320"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").Value = "'11.01"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +3.83%  "

$ws.Range("D14").Value = "'4.620.21"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").Value = "'3.991.42"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("E16").Value = "  +9.40%  "

$ws.Range("D17").Value = "'14.16"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").Value = "'20.46"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "'72.804.22"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("D21").Value = "'437.14"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +3.26%  "

$ws.Range("D22").Value = "'4.78"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +13.43%  "

$ws.Range("D23").Value = "'96.46"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("D24").Value = "'3.45"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -4.06%  "

$ws.Range("D25").Value = "'4.48"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +19.23%  "

$ws.Range("D26").Value = "'14.33"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").Value = "'11.32"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'5.96"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'10.53"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").Value = "'36.49"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "'7.96"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +2.33%  "

$ws.Range("D32").Value = "'13.74"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +2.82%  "

$ws.Range("D33").Value = "'0.132"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("D34").Value = "'48.80"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  -5.58%  "

$ws.Range("D35").Value = "'673.46"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("D36").Value = "'69.91"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +7.88%  "

$ws.Range("E37").Value = "  +7.70%  "

$ws.Range("D38").Value = "'0.439"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").Value = "'3.34"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +4.41%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'10.89"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +11.55%  "

$ws.Range("D45").Value = "'0.0490"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("E47").Value = "  -2.80%  "

$ws.Range("D48").Value = "'3.40"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").Value = "'3.03"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").Value = "'2.825.20"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +12.24%  "

$ws.Range("D51").Value = "'3.41"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +5.32%  "
